$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of portfolio data as row 85.
# Column A holds a date-looking string that must stay literal text
# (not get auto-converted to a date serial number), so we force the
# cell to Text format before assigning the value, then reset the
# style back to Normal so no stray formatting is left behind.
$row = 85
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-11-08"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = 57.38000106811523
$ws.Range("C$row").Value = 405.7000122070312
$ws.Range("D$row").Value = 306.1000061035156
